$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the four new "Database / SQL" related labels in column D, rows 17-20.
# These pull in four brand-new shared-string entries.
$ws.Range("D17").Value = " SQL "
$ws.Range("D18").Value = "Daily Activity - Per User"
$ws.Range("D19").Value = "All User Data - Secure"
$ws.Range("D20").Value = "Leave Data - Per User"

# Reposition the view: scroll so row 16 is at the top and select B26.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B26").Select()
